# edit.ps1 - apply methods.docx textual updates via Word COM Find/Replace
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Number of boat-based surveys: 15 -> 16
Replace-Text "combined data from 15 boat-based" "combined data from 16 boat-based"

# 2. discretewq package version bump (intro paragraph, with DOI-style citation kept identical)
Replace-Text "R package discretewq v2.3.2:" "R package discretewq v2.4.0:"

# 3. Number of fish-focused surveys: Eight -> Ten
Replace-Text "Davis. Eight surveys are primarily focused on collecting fish abundance" "Davis. Ten surveys are primarily focused on collecting fish abundance"

# 4. Key nutrient variables sentence gains "and other water quality parameters"
Replace-Text "Key nutrient variables were retained from the" "Key nutrient variables and other water quality parameters were retained from the"

# 5. Bottom water temperature -> water temperature and conductivity (was/were)
Replace-Text "In addition, water temperature from the bottom of the water column was retained when available." "In addition, water temperature and conductivity from the bottom of the water column were retained when available."

# 6. Bottom temperature samples -> bottom temperature and conductivity samples
Replace-Text "Bottom temperature samples were collected within 1 m of the bottom" "Bottom temperature and conductivity samples were collected within 1 m of the bottom"

# 7. SKT CDEC note: "they were transcribed" -> "they transcribed"
Replace-Text "SKT had notes on some temperature records that they were transcribed from a different monitoring program (CDEC)" "SKT had notes on some temperature records that they transcribed from a different monitoring program (CDEC)"

# 8. Microcystis surveys count: 3 -> 6
Replace-Text "by the 3 surveys that measured this variable" "by the 6 surveys that measured this variable"

# 9. Typo fix: "5-point scaled" -> "5-point scale"
Replace-Text "on the 5-point scaled, so all records" "on the 5-point scale, so all records"

# 10. Chlorophyll: EMP -> EMP and NCRO (filtered water samples)
Replace-Text "Chlorophyll-a methods differed slightly among surveys. EMP filtered water samples through a 1" "Chlorophyll-a methods differed slightly among surveys. EMP and NCRO filtered water samples through a 1"

# 11. Add comma: "in the field but USGS_SFBS" -> "in the field, but USGS_SFBS"
Replace-Text "used sonde probes to measure chlorophyll in the field but USGS_SFBS" "used sonde probes to measure chlorophyll in the field, but USGS_SFBS"

# 12. calibrated ... similar to EMP. -> similar to EMP and NCRO.
Replace-Text "collected and analyzed similar to EMP. " "collected and analyzed similar to EMP and NCRO. "

# 13. Nutrient protocol citation swap (EMP -> EMP and NCRO; Interagency... -> Battey and Perry 2023)
Replace-Text "EMP collected and preserved nutrients samples in accordance with standard protocols (Interagency Ecological Program et al. 2021a), after which they were processed in a lab." "EMP and NCRO collected and preserved nutrients samples in accordance with standard protocols (Battey and Perry 2023), after which they were processed in a lab."

# 14. Dissolved inorganic nutrients / pump sentence expansion (adds NCRO sampling method)
Replace-Text "collected, preserved, and processed dissolved inorganic nutrients in a similar manner to EMP. Both surveys collected water using a fixed flow-through pump." "collected, preserved, and processed dissolved inorganic nutrients in a similar manner to EMP and NCRO. EMP and USGS_SFBS collected water using a fixed flow-through pump, while NCRO collected water either by a bucket grab from the surface or a Van Dorn water sampler at a depth of approximately 1 meter."

# 15. discretewq package version bump (section 3, second citation of the R package)
Replace-Text "discretewq R package v2.3.2 (https://github.com/" "discretewq R package v2.4.0 (https://github.com/"

Write-Output "All replacements applied."
